$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 707 ("「諸行無常」..." post), which shifts all subsequent
# rows up by one (so old row 708 becomes new row 707, etc.)
$ws.Rows.Item(707).Delete()
